# New sensor data from sampling/scraping: append the latest reading row to
# each logger sheet, then restore each sheet's on-screen selection the way
# the author left it, finishing on the "ME" tab (the last sheet touched).

$wb = $excel.ActiveWorkbook

# --- PBSF: new row 37 ---
$ws = $wb.Worksheets.Item("PBSF")
$ws.Select()
$ws.Range("A37").Value = 44257.583333333336
$ws.Range("B37").Value = 1439
$ws.Range("C37").Value = 3.8
$ws.Range("C54:C55").Select()

# --- WIC: new row 23 ---
$ws = $wb.Worksheets.Item("WIC")
$ws.Select()
$ws.Range("A23").Value = 44257.440972222219
$ws.Range("B23").Value = 677.2
$ws.Range("C23").Value = 1.7
$ws.Range("K31").Select()

# --- YS: new row 38 ---
$ws = $wb.Worksheets.Item("YS")
$ws.Select()
$ws.Range("A38").Value = 44257.455555555556
$ws.Range("B38").Value = 340.2
$ws.Range("C38").Value = 2.1
$ws.Range("G36").Select()

# --- SW: new row 35 (datetime only, no conductivity/temp reading yet) ---
$ws = $wb.Worksheets.Item("SW")
$ws.Select()
$ws.Range("A35").Value = 44257.473611111112
$ws.Range("A35").Select()

# --- YI: new row 35 ---
$ws = $wb.Worksheets.Item("YI")
$ws.Select()
$ws.Range("A35").Value = 44257.489583333336
$ws.Range("B35").Value = 320.7
$ws.Range("C35").Value = 2.6
$ws.Range("C35").Select()

# --- YN: new row 34 ---
$ws = $wb.Worksheets.Item("YN")
$ws.Select()
$ws.Range("A34").Value = 44257.544444444444
$ws.Range("B34").Value = 461.8
$ws.Range("C34").Value = 1.6
$ws.Range("C39").Select()

# --- 6MC: new row 36 ---
$ws = $wb.Worksheets.Item("6MC")
$ws.Select()
$ws.Range("A36").Value = 44257.554861111108
$ws.Range("B36").Value = 481.6
$ws.Range("C36").Value = 3.7
$ws.Range("C36").Select()

# --- DC: new row 36 ---
$ws = $wb.Worksheets.Item("DC")
$ws.Select()
$ws.Range("A36").Value = 44257.560416666667
$ws.Range("B36").Value = 512.79999999999995
$ws.Range("C36").Value = 5.6
$ws.Range("E43").Select()

# --- PBMS: new row 37 ---
$ws = $wb.Worksheets.Item("PBMS")
$ws.Select()
$ws.Range("A37").Value = 44257.574999999997
$ws.Range("B37").Value = 870.2
$ws.Range("C37").Value = 3.1
$ws.Range("B41").Select()

# --- ME: no data change, just the last sheet the author left active ---
$ws = $wb.Worksheets.Item("ME")
$ws.Select()
$ws.Range("A12:A17").Select()
